$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain decimal number.
# Force them to remain text (matching the original inlineStr type in the OOXML)
# by flipping to Text format before the write, then clearing the format
# afterwards so no stray style index is left on the cell.
$textCells = @("D5", "D6", "D8", "D10", "D12", "D13", "D16", "D19", "D20", "D22", "D23", "D25", "D26", "D28", "D30", "D31", "D33", "D35", "D36", "D38", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "56.226.10"
$ws.Range("E2").Value = "  -3.12%  "
$ws.Range("D3").Value = "2.373.59"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "501.51"
$ws.Range("E5").Value = "  -4.83%  "
$ws.Range("D6").Value = "129.45"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "0.554"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").Value = "2.395.54"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").Value = "0.0958"
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.316"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").Value = "4.62"
$ws.Range("E13").Value = "  -7.04%  "
$ws.Range("D14").Value = "2.802.12"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "56.144.71"
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "21.48"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "2.428.46"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "10.09"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("D20").Value = "309.43"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  -3.61%  "
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").Value = "0.995"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("D25").Value = "64.27"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "2.488.65"
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("D28").Value = "0.374"
$ws.Range("E28").Value = "  -7.27%  "
$ws.Range("E29").Value = "  -5.21%  "
$ws.Range("D30").Value = "7.30"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "172.03"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "0.0₃0716"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "1.66"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "17.84"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").Value = "35.87"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D42").Value = "0.799"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.43"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "129.06"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("D45").Value = "3.34"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").Value = "4.75"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").Value = "250.83"
$ws.Range("E47").Value = "  -6.79%  "
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("D49").Value = "0.0903"
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("D50").Value = "0.0485"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "16.89"
$ws.Range("E51").Value = "  -1.08%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
